# Generate Report for Handoff
#
# - The handoff markdown file was regenerated under a new GUID name and its
#   status flipped from "Ready for handoff" to "Handoff transform failed".
# - Because the transform failed, the per-language rows no longer have a
#   "Latest Handoff File" / "Latest Handoff Datetime" result: the handoff
#   file cell + its hyperlink are cleared, the handoff datetime resets to
#   the zero-date, and the per-row status flips from "Include" to "Ignored".

function Remove-HyperlinkAt($ws, $addr) {
    $target = $null
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $target = $hl
        }
    }
    if ($target -ne $null) {
        $target.Delete()
    }
}

function Set-HyperlinkDisplay($ws, $addr, $text) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $hl.TextToDisplay = $text
        }
    }
}

$wb = $excel.ActiveWorkbook

$newName = "88e15722-ee89-4119-8452-37afc5243a5f.md"
$newStatus = "Handoff transform failed"
$zeroDate = "0001-01-01 00:00:00"

# --- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newName
Set-HyperlinkDisplay $wsOverview '$A$2' $newName
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus

# --- zh-cn sheet --------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newName
Set-HyperlinkDisplay $wsZh '$A$2' $newName
$wsZh.Range("B2").Value = $newStatus

Remove-HyperlinkAt $wsZh '$C$2'
$wsZh.Range("C2").Clear()

$wsZh.Range("D2").Value = $zeroDate
$wsZh.Range("H2").Value = "Ignored"

# --- de-de sheet --------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newName
Set-HyperlinkDisplay $wsDe '$A$2' $newName
$wsDe.Range("B2").Value = $newStatus

Remove-HyperlinkAt $wsDe '$C$2'
$wsDe.Range("C2").Clear()

$wsDe.Range("D2").Value = $zeroDate
$wsDe.Range("H2").Value = "Ignored"
